$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.50699654812515
$ws.Range("C2").Value = 10.64534155274959
$ws.Range("D2").Value = 12.2463787708175
$ws.Range("F2").Value = 26.19612226925073
$ws.Range("G2").Value = 21.77841439143124
$ws.Range("H2").Value = 12.57751115334616
$ws.Range("I2").Value = 16.92113940241933
$ws.Range("J2").Value = 11.25502089938484
$ws.Range("O2").Value = 18.20247603018952
$ws.Range("B3").Value = 14.68668325018348
$ws.Range("C3").Value = 9.985671679390919
$ws.Range("D3").Value = 12.13780463539512
$ws.Range("F3").Value = 26.26434819624405
$ws.Range("G3").Value = 21.90666694744682
$ws.Range("H3").Value = 12.64318489558068
$ws.Range("I3").Value = 17.0798041765912
$ws.Range("J3").Value = 11.23424606417158
$ws.Range("O3").Value = 18.31531901186841
$ws.Range("B4").Value = 14.15775925236431
$ws.Range("C4").Value = 9.556117943716934
$ws.Range("D4").Value = 12.07271666721393
$ws.Range("F4").Value = 26.31529312649378
$ws.Range("G4").Value = 21.9972863988707
$ws.Range("H4").Value = 12.68623187677778
$ws.Range("I4").Value = 17.18207778689885
$ws.Range("J4").Value = 11.22376933310073
$ws.Range("O4").Value = 18.39029633780381
$ws.Range("B5").Value = 13.93604057776153
$ws.Range("C5").Value = 9.374934966272598
$ws.Range("D5").Value = 12.04661370161252
$ws.Range("F5").Value = 26.33831939053733
$ws.Range("G5").Value = 22.03716937751358
$ws.Range("H5").Value = 12.70445809101561
$ws.Range("I5").Value = 17.22497851360752
$ws.Range("J5").Value = 11.22007602236737
$ws.Range("O5").Value = 18.42227563440161
$ws.Range("B6").Value = 13.89885690138701
$ws.Range("C6").Value = 9.34448006832239
$ws.Range("D6").Value = 12.04230544286845
$ws.Range("F6").Value = 26.34227939474839
$ws.Range("G6").Value = 22.04396946902625
$ws.Range("H6").Value = 12.70752585091462
$ws.Range("I6").Value = 17.23217612915019
$ws.Range("J6").Value = 11.21949761918819
$ws.Range("O6").Value = 18.42767169915272
$ws.Range("B7").Value = 14.15479383995668
$ws.Range("C7").Value = 9.553699248393436
$ws.Range("D7").Value = 12.07236289762634
$ws.Range("F7").Value = 26.31559450779322
$ws.Range("G7").Value = 21.99781235177727
$ws.Range("H7").Value = 12.68647491215915
$ws.Range("I7").Value = 17.18265140345921
$ws.Range("J7").Value = 11.22371718791069
$ws.Range("O7").Value = 18.39072185778598
$ws.Range("B8").Value = 15.22948540806203
$ws.Range("C8").Value = 10.42298888499467
$ws.Range("D8").Value = 12.20862994319572
$ws.Range("F8").Value = 26.21776117065962
$ws.Range("G8").Value = 21.82015387799969
$ws.Range("H8").Value = 12.59958998158739
$ws.Range("I8").Value = 16.97484170415712
$ws.Range("J8").Value = 11.24738661318213
$ws.Range("O8").Value = 18.24019963339936
$ws.Range("B9").Value = 17.13065215773294
$ws.Range("C9").Value = 11.93240106507209
$ws.Range("D9").Value = 12.48725945304468
$ws.Range("F9").Value = 26.09819034675818
$ws.Range("G9").Value = 21.56726849384085
$ws.Range("H9").Value = 12.45083805759301
$ws.Range("I9").Value = 16.60568806709608
$ws.Range("J9").Value = 11.31173680942858
$ws.Range("O9").Value = 17.99043953607497
$ws.Range("B10").Value = 18.39546900999768
$ws.Range("C10").Value = 12.92182824683535
$ws.Range("D10").Value = 12.69741212238864
$ws.Range("F10").Value = 26.0549385068845
$ws.Range("G10").Value = 21.44138299805644
$ws.Range("H10").Value = 12.35476784586243
$ws.Range("I10").Value = 16.35765246224821
$ws.Range("J10").Value = 11.36971993563802
$ws.Range("O10").Value = 17.83497818118336
$ws.Range("B11").Value = 18.94131213472965
$ws.Range("C11").Value = 13.34598911536303
$ws.Range("D11").Value = 12.79386949925393
$ws.Range("F11").Value = 26.04503553008625
$ws.Range("G11").Value = 21.39745207105339
$ws.Range("H11").Value = 12.31394081608832
$ws.Range("I11").Value = 16.24980517677766
$ws.Range("J11").Value = 11.39836420093192
$ws.Range("O11").Value = 17.77041942392786
$ws.Range("B12").Value = 19.14370786771673
$ws.Range("C12").Value = 13.50288236226986
$ws.Range("D12").Value = 12.83049098541823
$ws.Range("F12").Value = 26.04269610767441
$ws.Range("G12").Value = 21.38275847440698
$ws.Range("H12").Value = 12.29889500467186
$ws.Range("I12").Value = 16.20967996927637
$ws.Range("J12").Value = 11.40953142060254
$ws.Range("O12").Value = 17.74686518067351
$ws.Range("B13").Value = 19.10031048563166
$ws.Range("C13").Value = 13.46925830234145
$ws.Range("D13").Value = 12.82260010816901
$ws.Range("F13").Value = 26.04313713558075
$ws.Range("G13").Value = 21.38583625940874
$ws.Range("H13").Value = 12.30211693864508
$ws.Range("I13").Value = 16.21828992753261
$ws.Range("J13").Value = 11.40711221341912
$ws.Range("O13").Value = 17.75189819484274
$ws.Range("B14").Value = 18.95804987782726
$ws.Range("C14").Value = 13.3589714951784
$ws.Range("D14").Value = 12.79688062472332
$ws.Range("F14").Value = 26.04481476761044
$ws.Range("G14").Value = 21.39620416432481
$ws.Range("H14").Value = 12.31269467529344
$ws.Range("I14").Value = 16.24648975290878
$ws.Range("J14").Value = 11.39927655112325
$ws.Range("O14").Value = 17.76846367115867
$ws.Range("B15").Value = 18.87034926806214
$ws.Range("C15").Value = 13.29093235175455
$ws.Range("D15").Value = 12.78113825742877
$ws.Range("F15").Value = 26.04602620586243
$ws.Range("G15").Value = 21.40280841657928
$ws.Range("H15").Value = 12.31922785367071
$ws.Range("I15").Value = 16.26385589422219
$ws.Range("J15").Value = 11.39451850925077
$ws.Range("O15").Value = 17.77872696814267
$ws.Range("B16").Value = 18.35919751381395
$ws.Range("C16").Value = 12.89358657852718
$ws.Range("D16").Value = 12.69112325719976
$ws.Range("F16").Value = 26.05578282485056
$ws.Range("G16").Value = 21.44452446507731
$ws.Range("H16").Value = 12.35749393797023
$ws.Range("I16").Value = 16.3648006247843
$ws.Range("J16").Value = 11.36789310676496
$ws.Range("O16").Value = 17.83932180749447
$ws.Range("B17").Value = 18.03801469591115
$ws.Range("C17").Value = 12.64318588971494
$ws.Range("D17").Value = 12.63610161785107
$ws.Range("F17").Value = 26.06427539206682
$ws.Range("G17").Value = 21.47354941404961
$ws.Range("H17").Value = 12.38170617296845
$ws.Range("I17").Value = 16.42800190171332
$ws.Range("J17").Value = 11.35213604734097
$ws.Range("O17").Value = 17.87807765476636
$ws.Range("B18").Value = 17.8505019716613
$ws.Range("C18").Value = 12.4967206197522
$ws.Range("D18").Value = 12.6045370603571
$ws.Range("F18").Value = 26.07007968962965
$ws.Range("G18").Value = 21.49149760585716
$ws.Range("H18").Value = 12.39590303303506
$ws.Range("I18").Value = 16.46482300623667
$ws.Range("J18").Value = 11.34328680210201
$ws.Range("O18").Value = 17.90094852590447
$ws.Range("B19").Value = 17.78653841175577
$ws.Range("C19").Value = 12.44671066569723
$ws.Range("D19").Value = 12.59386486164566
$ws.Range("F19").Value = 26.07220268218027
$ws.Range("G19").Value = 21.49778911790501
$ws.Range("H19").Value = 12.40075629699139
$ws.Range("I19").Value = 16.47737070191091
$ws.Range("J19").Value = 11.34032748930855
$ws.Range("O19").Value = 17.90879153998939
$ws.Range("B20").Value = 18.07249300441334
$ws.Range("C20").Value = 12.67009414196964
$ws.Range("D20").Value = 12.64195043810081
$ws.Range("F20").Value = 26.06327612591104
$ws.Range("G20").Value = 21.47032971003756
$ws.Range("H20").Value = 12.37910072390362
$ws.Range("I20").Value = 16.42122545604828
$ws.Range("J20").Value = 11.35379132648137
$ws.Range("O20").Value = 17.87389200033929
$ws.Range("B21").Value = 18.99995244628752
$ws.Range("C21").Value = 13.3914665300966
$ws.Range("D21").Value = 12.80443270173507
$ws.Range("F21").Value = 26.04428368662811
$ws.Range("G21").Value = 21.39310596610299
$ws.Range("H21").Value = 12.30957648037458
$ws.Range("I21").Value = 16.23818741794889
$ws.Range("J21").Value = 11.40156943173159
$ws.Range("O21").Value = 17.76357370602704
$ws.Range("B22").Value = 19.58100048359967
$ws.Range("C22").Value = 13.84120257206384
$ws.Range("D22").Value = 12.91116600270492
$ws.Range("F22").Value = 26.04009472918874
$ws.Range("G22").Value = 21.35396568991917
$ws.Range("H22").Value = 12.26655485988964
$ws.Range("I22").Value = 16.12272329057875
$ws.Range("J22").Value = 11.434658902457
$ws.Range("O22").Value = 17.69668115192311
$ws.Range("B23").Value = 19.27319654356414
$ws.Range("C23").Value = 13.60315619512986
$ws.Range("D23").Value = 12.85416014318075
$ws.Range("F23").Value = 26.04157656717802
$ws.Range("G23").Value = 21.37381133463634
$ws.Range("H23").Value = 12.28929487035536
$ws.Range("I23").Value = 16.18396873125217
$ws.Range("J23").Value = 11.41682992585995
$ws.Range("O23").Value = 17.73190430113004
$ws.Range("B24").Value = 18.05691426634444
$ws.Range("C24").Value = 12.6579367287471
$ws.Range("D24").Value = 12.63930597214576
$ws.Range("F24").Value = 26.06372502317175
$ws.Range("G24").Value = 21.47178141103363
$ws.Range("H24").Value = 12.3802777848783
$ws.Range("I24").Value = 16.42428757350364
$ws.Range("J24").Value = 11.35304232116315
$ws.Range("O24").Value = 17.87578249856351
$ws.Range("B25").Value = 16.63914814481544
$ws.Range("C25").Value = 11.54504315196113
$ws.Range("D25").Value = 12.41081580710405
$ws.Range("F25").Value = 26.12273975698203
$ws.Range("G25").Value = 21.62527644894307
$ws.Range("H25").Value = 12.48876067588677
$ws.Range("I25").Value = 16.70146840321311
$ws.Range("J25").Value = 11.29243005673429
$ws.Range("O25").Value = 18.05310827515256
